$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3   = -13.224
    14  = -12.139
    16  = -12.983
    21  = -12.842
    23  = -12.223
    25  = -12.682
    26  = -12.576
    29  = -11.941
    40  = -12.426
    53  = -11.443
    57  = -13.829
    59  = -13.155
    65  = -12.118
    69  = -10.676
    79  = -12.083
    83  = -13.3
    91  = -10.595
    93  = -11.764
    100 = -13.395
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
